$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh scraped crypto data: Price (D) and 1h Volume change (E) columns.
# Price cells are stored as literal text in the source sheet (mixed locale
# formatting like "2.483.83"), so force text format before writing each one
# to stop Excel from re-interpreting the string as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.089.13"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.482.48"
$ws.Range("E3").Value = "  +2.81%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.38"
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("E6").Value = "  +4.01%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +2.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.483.05"
$ws.Range("E9").Value = "  +2.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  +5.75%  "

$ws.Range("E11").Value = "  +1.81%  "

$ws.Range("E12").Value = "  +4.05%  "

$ws.Range("E13").Value = "  +3.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.42"
$ws.Range("E14").Value = "  +3.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.887.72"
$ws.Range("E15").Value = "  +5.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.971.92"
$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("E17").Value = "  +3.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.388.42"
$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.02"
$ws.Range("E19").Value = "  +0.47%  "

$ws.Range("E20").Value = "  +1.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.89"
$ws.Range("E21").Value = "  +0.56%  "

$ws.Range("E22").Value = "  +2.31%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("E24").Value = "  +3.03%  "

$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("E26").Value = "  +4.50%  "

$ws.Range("E27").Value = "  +5.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.622.39"
$ws.Range("E29").Value = "  +4.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0907"
$ws.Range("E30").Value = "  +4.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "511.60"
$ws.Range("E31").Value = "  +5.46%  "

$ws.Range("E32").Value = "  +1.39%  "

$ws.Range("E33").Value = "  +3.63%  "

$ws.Range("E34").Value = "  +0.90%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.33"
$ws.Range("E36").Value = "  +2.67%  "

$ws.Range("E37").Value = "  +5.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.70"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.24"
$ws.Range("E39").Value = "  +0.89%  "

$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("E41").Value = "  +3.59%  "

$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("E43").Value = "  +2.98%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.80"
$ws.Range("E44").Value = "  +5.17%  "

$ws.Range("E45").Value = "  +3.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.86"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.79"
$ws.Range("E47").Value = "  +4.19%  "

$ws.Range("E48").Value = "  +1.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.515"
$ws.Range("E49").Value = "  +2.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0251"
$ws.Range("E50").Value = "  +4.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0728"
$ws.Range("E51").Value = "  +1.60%  "
